# Generate Report for Handoff
# Updates the localization-status report: the "In Translation" rows have
# moved to "Ready for handoff", and the handoff / generation timestamps
# are refreshed. Column widths on the affected "status" columns are
# widened to fit the new (longer) status text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-15 09:12:48"
$overview.Range("G2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$overview.Columns.Item(5).ColumnWidth = 16.33
$overview.Columns.Item(6).ColumnWidth = 16.33

# --- zh-cn sheet --------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-15 09:12:44"
$zhcn.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zhcn.Columns.Item(3).ColumnWidth = 16.33

# --- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-15 09:12:48"
$dede.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$dede.Columns.Item(3).ColumnWidth = 16.33
